# Update the yearly income-statement "database" (rolling the 5-year window
# forward by one fiscal year) and refresh the published-on dates, per the
# commit "update database and change read_price algorithm".
#
# Columns D:H hold five consecutive fiscal years. Every column's header
# (row 8), publish date (row 9) and figures (rows 11-27) shift to what used
# to be in the next column, and column H receives the brand-new year's data.
# The read_price algorithm change is reflected by the new figures now being
# expressed in (small, e.g. thousand/rial-per-dollar) units rather than the
# previous large rial amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: fiscal period headers --------------------------------------
$ws.Cells.Item(8, 4).Value = "12 ماهه منتهی به 1397/12"
$ws.Cells.Item(8, 5).Value = "12 ماهه منتهی به 1398/12"
$ws.Cells.Item(8, 6).Value = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 7).Value = "12 ماهه منتهی به 1400/12"
$ws.Cells.Item(8, 8).Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ------------------------------------------------
$ws.Cells.Item(9, 4).Value = "1399-03-07 (8)"
$ws.Cells.Item(9, 5).Value = "1400-02-28 (9)"
$ws.Cells.Item(9, 6).Value = "1401-03-11 (8)"
$ws.Cells.Item(9, 7).Value = "1402-02-30 (7)"
$ws.Cells.Item(9, 8).Value = "1402-02-30"

# --- Row 11: فروش (Sales) -------------------------------------------------
$ws.Cells.Item(11, 4).Value = 26625
$ws.Cells.Item(11, 5).Value = 31425
$ws.Cells.Item(11, 6).Value = 37454
$ws.Cells.Item(11, 7).Value = 51890
$ws.Cells.Item(11, 8).Value = 69840

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ----------
$ws.Cells.Item(12, 4).Value = -26583
$ws.Cells.Item(12, 5).Value = -28028
$ws.Cells.Item(12, 6).Value = -32346
$ws.Cells.Item(12, 7).Value = -42900
$ws.Cells.Item(12, 8).Value = -56066

# --- Row 13: سود (زیان) ناخالص (Gross profit) -----------------------------
$ws.Cells.Item(13, 4).Value = "-"
$ws.Cells.Item(13, 5).Value = 3397
$ws.Cells.Item(13, 6).Value = 5108
$ws.Cells.Item(13, 7).Value = 8990
$ws.Cells.Item(13, 8).Value = 13774

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ------------
$ws.Cells.Item(14, 4).Value = -2483
$ws.Cells.Item(14, 5).Value = -2761
$ws.Cells.Item(14, 6).Value = -3503
$ws.Cells.Item(14, 7).Value = -4144
$ws.Cells.Item(14, 8).Value = -3853

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (هزینه استثنایی) ------------------
$ws.Cells.Item(15, 4).Value = "-"
$ws.Cells.Item(15, 5).Value = "-"
$ws.Cells.Item(15, 6).Value = "-"
$ws.Cells.Item(15, 7).Value = "-"
$ws.Cells.Item(15, 8).Value = "-"

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ------------------------
$ws.Cells.Item(16, 4).Value = 4047
$ws.Cells.Item(16, 5).Value = 1882
$ws.Cells.Item(16, 6).Value = 1820
$ws.Cells.Item(16, 7).Value = 17
$ws.Cells.Item(16, 8).Value = 2503

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ------------------------
$ws.Cells.Item(17, 4).Value = "-"
$ws.Cells.Item(17, 5).Value = 2518
$ws.Cells.Item(17, 6).Value = 3425
$ws.Cells.Item(17, 7).Value = 4864
$ws.Cells.Item(17, 8).Value = 12423

# --- Row 18: هزینه های مالی (Financial expenses) --------------------------
$ws.Cells.Item(18, 4).Value = -46
$ws.Cells.Item(18, 5).Value = -34
$ws.Cells.Item(18, 6).Value = -76
$ws.Cells.Item(18, 7).Value = -183
$ws.Cells.Item(18, 8).Value = -44

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ----------------------
$ws.Cells.Item(19, 4).Value = 119
$ws.Cells.Item(19, 5).Value = 183
$ws.Cells.Item(19, 6).Value = 659
$ws.Cells.Item(19, 7).Value = 817
$ws.Cells.Item(19, 8).Value = 400

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات -------------
$ws.Cells.Item(20, 4).Value = "-"
$ws.Cells.Item(20, 5).Value = 2667
$ws.Cells.Item(20, 6).Value = 4008
$ws.Cells.Item(20, 7).Value = 5498
$ws.Cells.Item(20, 8).Value = 12779

# --- Row 21: مالیات (Tax) --------------------------------------------------
$ws.Cells.Item(21, 4).Value = "-"
$ws.Cells.Item(21, 5).Value = "-"
$ws.Cells.Item(21, 6).Value = -22
$ws.Cells.Item(21, 7).Value = -352
$ws.Cells.Item(21, 8).Value = -1448

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---------------------------
$ws.Cells.Item(22, 4).Value = "-"
$ws.Cells.Item(22, 5).Value = 2667
$ws.Cells.Item(22, 6).Value = 3985
$ws.Cells.Item(22, 7).Value = 5146
$ws.Cells.Item(22, 8).Value = 11331

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی -----------------
$ws.Cells.Item(23, 4).Value = "-"
$ws.Cells.Item(23, 5).Value = "-"
$ws.Cells.Item(23, 6).Value = "-"
$ws.Cells.Item(23, 7).Value = "-"
$ws.Cells.Item(23, 8).Value = "-"

# --- Row 24: سود (زیان) خالص (Net profit) ----------------------------------
$ws.Cells.Item(24, 4).Value = "-"
$ws.Cells.Item(24, 5).Value = 2667
$ws.Cells.Item(24, 6).Value = 3985
$ws.Cells.Item(24, 7).Value = 5146
$ws.Cells.Item(24, 8).Value = 11331

# --- Row 25: سود هر سهم پس از کسر مالیات (EPS) -----------------------------
$ws.Cells.Item(25, 4).Value = "-"
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0

# --- Row 26: سرمایه (Capital) ----------------------------------------------
$ws.Cells.Item(26, 4).Value = 2106
$ws.Cells.Item(26, 5).Value = 4389
$ws.Cells.Item(26, 6).Value = 8076
$ws.Cells.Item(26, 7).Value = 6921
$ws.Cells.Item(26, 8).Value = 5174

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه -------------------------------
$ws.Cells.Item(27, 4).Value = "-"
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 0
